# Dispatcher completed and tested
# Update the "Loan Data" sheet: rename the income column header, refresh the
# sample data rows, drop the centered cell styling back to the default
# "Normal" style, narrow column C, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: "Current Yearly Income" -> "Yearly Income" ---
$ws.Range("C1").Value = "Yearly Income"

# --- Refresh sample data rows (A: Loan Amount Requested, B: Loan Term,
#     C: Yearly Income, D: Age) ---
$ws.Range("A2").Value = 50000
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 10000
$ws.Range("D2").Value = 23

$ws.Range("A3").Value = 3000
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 20000
$ws.Range("D3").Value = 65

$ws.Range("A4").Value = 450000
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 700000
$ws.Range("D4").Value = 33

$ws.Range("A5").Value = 90000
$ws.Range("B5").Value = 7
$ws.Range("C5").Value = 130000
$ws.Range("D5").Value = 45

# --- Reset cell styling on the used range back to the default "Normal"
#     style, removing the centered-alignment formatting ---
$ws.Range("A1:D5").Style = "Normal"

# --- Narrow column C to fit the shorter header/values ---
$ws.Columns("C").ColumnWidth = 12.67

# --- Move the active selection to F4 ---
$ws.Range("F4").Select()
